# Generate Report for Handoff
#
# Source file "b" (zh-cn / de-de) has been handed off again: its status
# moves from "Handed back: in sync with en-US" to "Ready for handoff", a
# new handoff file + timestamp is recorded, and the "Latest Handoff File"
# hyperlink's visible label is updated to point at the new handoff file
# name (the underlying hyperlink target URL is left exactly as-is).

$wb = $excel.ActiveWorkbook

$newStatus        = "Ready for handoff"
$zhHandoffFile    = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhHandoffDate    = "2016-02-23 08:50:27"
$deHandoffFile    = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$deHandoffDate    = "2016-02-23 08:50:38"

# ---------------------------------------------------------------------
# 1) Overview sheet: b.md row's per-locale status (B3 = zh-cn, C3 = de-de)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2) zh-cn sheet: b.md row (row 3) handoff info
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B3").Value = $newStatus
$zh.Range("C3").Value = $zhHandoffFile
$zh.Range("D3").Value = $zhHandoffDate

# Rebuild the sheet's hyperlinks so the "Latest Handoff File" link on C3
# shows the new file name. Item-level hyperlink edits in this host always
# append a new link object rather than amending one in place, and the only
# way to truly drop a stale entry is to clear the whole collection - so the
# other (unchanged) links are recreated with their original address/text.
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2e4df7b2a22f28cfaff0cdc755cecbaadfb9b3fd/e2e/a.md", "", "", "a.md")
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9f14ae3be5ed4e70f0f19c6ac1bfb67c49a3a9ba/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/573868ca459fd995d7468c77b210228312040595/e2e/a.md", "", "", "a.md")
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/784a919c80a0ccc6ba45ed16e311232ccf4684d0/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2e4df7b2a22f28cfaff0cdc755cecbaadfb9b3fd/e2e/b.md", "", "", "b.md")
$zh.Hyperlinks.Add($zh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9f14ae3be5ed4e70f0f19c6ac1bfb67c49a3a9ba/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", $zhHandoffFile)
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/573868ca459fd995d7468c77b210228312040595/e2e/a.md", "", "", "a.md")
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/784a919c80a0ccc6ba45ed16e311232ccf4684d0/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/2e4df7b2a22f28cfaff0cdc755cecbaadfb9b3fd/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# 3) de-de sheet: b.md row (row 3) handoff info
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("B3").Value = $newStatus
$de.Range("C3").Value = $deHandoffFile
$de.Range("D3").Value = $deHandoffDate

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2e4df7b2a22f28cfaff0cdc755cecbaadfb9b3fd/e2e/a.md", "", "", "a.md")
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a370d7c138e72ff3091517048e980fa54bb17dbf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/a346b16d78538d7c36ad4392346fd334f0b8b94f/e2e/a.md", "", "", "a.md")
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b58cbc7c609d6cb59877eccd184dca8c62c820d7/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2e4df7b2a22f28cfaff0cdc755cecbaadfb9b3fd/e2e/b.md", "", "", "b.md")
$de.Hyperlinks.Add($de.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a370d7c138e72ff3091517048e980fa54bb17dbf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", $deHandoffFile)
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/a346b16d78538d7c36ad4392346fd334f0b8b94f/e2e/a.md", "", "", "a.md")
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b58cbc7c609d6cb59877eccd184dca8c62c820d7/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/2e4df7b2a22f28cfaff0cdc755cecbaadfb9b3fd/.localization-config", "", "", ".localization-config")
